$d = $word.ActiveDocument

# Step 1: trim the trailing period off the existing sentence so the new
# continuation can be appended seamlessly (matches the diff's first hunk).
$find = $d.Content.Find
$find.Execute(" __, when I probably should have just got started sooner.", $true, $false, $false, $false, $false, $true, 1, $false, " __, when I probably should have just got started sooner", 2)

# Step 2: locate the (now period-less) end of that sentence so we can append
# the new runs right after it.
$find2 = $d.Content.Find
$find2.Execute("just got started sooner", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchor = $d.Range($find2.Parent.End, $find2.Parent.End)
$anchor.Collapse(0)

$anchor.InsertAfter(", especially since I ran into a lot of trouble spots that __")
$anchor.Collapse(0)

$anchor.InsertAfter(".")
$anchor.Collapse(0)

$anchor.InsertAfter(" For example, in the activity-specific adjustments section of the fitness function, I ran into issues with ensuring that I was checking both instances of the two-section course against the two instances of the ")
$anchor.Collapse(0)

$italicStart = $anchor.End
$anchor.InsertAfter("other")
$italicEnd = $anchor.End
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Italic = $true
$anchor.Collapse(0)

$anchor.InsertAfter(" two-section course. ")
$anchor.Collapse(0)

$anchor.InsertAfter("I couldn’t afford to spend ")
$anchor.Collapse(0)

$anchor.InsertAfter("too much ")
$anchor.Collapse(0)

$anchor.InsertAfter("time")
$anchor.Collapse(0)

$anchor.InsertAfter(" on this issue, so eventually I just went with a simplistic method that didn’t entirely fulfill the specifications. ")
$anchor.Collapse(0)
